$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New row 23 content
$ws.Range("A23").Value2 = "body"
$ws.Range("B23").Value2 = 11

# 2. Copy A11's current format (style idx 6: sz16/theme1 font, fillId=8) onto A23
$ws.Range("A11").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats

# 3. Make A23 italic + size 12 (default body text size)
$ws.Range("A23").Font.Size = 12
$ws.Range("A23").Font.Italic = $true

# 4. Change A11's font color theme to Light2 (ooxml theme=2)
$ws.Range("A11").Font.ThemeColor = 4

# 5. Move selection to A11
$ws.Range("A11").Select()
